$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Table spacing tweaks: every table in the document gets its left indent
#    nudged from 47 -> 45 dxa (2.35pt -> 2.25pt) and its left cell margin
#    (table default + every individual cell) nudged from 45 -> 42 dxa
#    (2.25pt -> 2.1pt).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables($i)
    $t.Rows.LeftIndent = 2.25
    $t.LeftPadding = 2.1
    foreach ($row in $t.Rows) {
        foreach ($cell in $row.Cells) {
            $cell.LeftPadding = 2.1
        }
    }
}

# ---------------------------------------------------------------------------
# 2. The big course-schedule table (3rd table) gets its first / last column
#    widths nudged by 1 dxa: 447 -> 446 and 3051 -> 3052 (0.05pt each way).
# ---------------------------------------------------------------------------
$schedule = $d.Tables(3)
$schedule.Columns(1).Width = 22.3
$schedule.Columns(4).Width = 152.6

# ---------------------------------------------------------------------------
# 3. Fill in the readings / homework for item 5 "Clause Types" row.
# ---------------------------------------------------------------------------
$schedule.Cell(7, 3).Range.Text = "SIEG Ch.9 sections 9.1 -9.4"
$schedule.Cell(7, 4).Range.Text = "exercises 1, 2, 5 pg 173"

# ---------------------------------------------------------------------------
# 4. Normal style: allow punctuation to extend past the text margin
#    (w:overflowPunct false -> true).
# ---------------------------------------------------------------------------
$d.Styles("Normal").ParagraphFormat.HangingPunctuation = $true
